# Daily attendance processing - 2025-10-20 17:18:17
# Reconciles the Y3 CNS session-analysis sheet: re-orders the "Recorded By"
# attendee lists (list membership unchanged, order refreshed to match the
# latest export), refreshes a handful of roll-up counters/percentages, marks
# the Year 3 / C2 / HISTOLOGY session 2 (row 42) as now Recorded with its
# attendance + recorder, and narrows column I back to a width of 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write literal text into a cell WITHOUT letting the "looks like a
# number/percentage" auto-detection turn it into a numeric value, and
# WITHOUT disturbing the cell's existing style (fill/font/alignment).
#
# Approach: build the literal string as a formula result (string
# concatenation always yields a text value, never auto-converted) in an
# always-empty scratch cell (column J is a blank spacer column in this
# sheet), copy *values only* into the destination (so the destination's
# own formatting/style is left completely untouched), then fully clear
# the scratch cell again.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)

    $scratch = $ws.Range("J1")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# "Recorded By" attendee-list re-orderings (same attendees, new order)
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G5").Value = "hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
$ws.Range("G6").Value = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

$ws.Range("G12").Value = "salma.elgendy.std@med.asu.edu.eg, System"
$ws.Range("G13").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G24").Value = "yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G25").Value = "abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg"

$ws.Range("G30").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G31").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G32").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G33").Value = "hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
$ws.Range("G34").Value = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

$ws.Range("G40").Value = "salma.elgendy.std@med.asu.edu.eg, System"
$ws.Range("G41").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G52").Value = "yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G53").Value = "abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg"

# ---------------------------------------------------------------------
# H14: attendance fraction text ("56/221" -> "57/221") - not number-like,
# plain assignment keeps it text.
# ---------------------------------------------------------------------
$ws.Range("H14").Value = "57/221"

# ---------------------------------------------------------------------
# Numeric roll-up counters (plain numbers, no auto-detect risk)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 0
$ws.Range("O16").Value = 11
$ws.Range("P16").Value = 0

# ---------------------------------------------------------------------
# Percentage cells are stored as literal text (e.g. "39.3%"), not numbers.
# A direct .Value assignment of a percent-looking string gets silently
# reinterpreted as a numeric percentage, so route these through the
# text-preserving helper.
# ---------------------------------------------------------------------
Set-LiteralText $ws.Range("L9") "39.3%"
Set-LiteralText $ws.Range("L10") "42.3%"
Set-LiteralText $ws.Range("R16") "39.3%"
Set-LiteralText $ws.Range("S16") "37.0%"

# ---------------------------------------------------------------------
# Row 42 (Year 3 / C2 / HISTOLOGY session 2): was "Pending"-styled with no
# recorder; now Recorded. Re-stripe it with the same fill/font used by the
# other "Recorded" rows (copy formats only from row 41) and fill in the
# recorder / attendance / status.
# ---------------------------------------------------------------------
$ws.Range("A41:I41").Copy()
$ws.Range("A42:I42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G42").Value = "Safa.hany@med.asu.edu.eg"
$ws.Range("H42").Value = "66/246"
$ws.Range("I42").Value = "Recorded"

# ---------------------------------------------------------------------
# Column I width: 14 -> 10 chars. Copy column H's exact ColumnWidth (it is
# already the target width) instead of hard-coding a float, so the stored
# <col width="..."> comes out exactly "10" like the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth
